# Update "人气" (popularity/F column) figures on gh-pages regeneration.
# Mirrors the commit "Update gh-pages to output generated at 456a3b4":
# the F-column numeric values increased slightly on the 展览 (Exhibition),
# 演出 (Performance) and 全部类型 (All types) sheets.

$wb = $excel.ActiveWorkbook

# --- 展览 (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 1354
$ws1.Range("F8").Value  = 189
$ws1.Range("F10").Value = 8418
$ws1.Range("F14").Value = 280
$ws1.Range("F15").Value = 309
$ws1.Range("F18").Value = 345
$ws1.Range("F19").Value = 10604
$ws1.Range("F21").Value = 282
$ws1.Range("F35").Value = 881
$ws1.Range("F37").Value = 274
$ws1.Range("F42").Value = 752
$ws1.Range("F44").Value = 328
$ws1.Range("F45").Value = 283

# --- 演出 (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F10").Value = 29
$ws2.Range("F23").Value = 24

# --- 全部类型 (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 1354
$ws4.Range("F9").Value  = 189
$ws4.Range("F11").Value = 8418
$ws4.Range("F13").Value = 29
$ws4.Range("F16").Value = 280
$ws4.Range("F17").Value = 309
$ws4.Range("F19").Value = 345
$ws4.Range("F20").Value = 10604
$ws4.Range("F22").Value = 282
$ws4.Range("F34").Value = 881
$ws4.Range("F36").Value = 274
$ws4.Range("F41").Value = 752
$ws4.Range("F43").Value = 328
$ws4.Range("F44").Value = 24
$ws4.Range("F45").Value = 283

$wb.Save()
